# Generate Report for Handback
# Update the "generated" timestamps recorded on the handback status report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview!G2 - "Latest HO Xliff Generate Date" for the first file.
# This shared string is also used by de-de!H2, so both cells move together.
$overview.Range("G2").Value = "2016-08-31 05:07:20"

# zh-cn sheet: Correspond Handoff/Handback datetimes for the first row.
$zhcn.Range("H2").Value = "2016-08-31 05:07:15"
$zhcn.Range("K2").Value = "2016-08-31 05:07:33"

# de-de sheet: Correspond Handoff datetime shares text with Overview!G2 (updated above),
# and the Correspond Handback datetime is updated independently.
$dede.Range("H2").Value = "2016-08-31 05:07:20"
$dede.Range("K2").Value = "2016-08-31 05:07:40"
